$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78-117 down to 79-118
$ws.Rows.Item(78).EntireRow.Insert()

# Populate the newly inserted row 78 with the new weekly price record
$ws.Cells.Item(78, 1).Value = 7
$ws.Cells.Item(78, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(78, 3).Value = "Ñuble"
$ws.Cells.Item(78, 4).Value = 44957
$ws.Cells.Item(78, 5).Value = 16
$ws.Cells.Item(78, 6).Value = 100112031
$ws.Cells.Item(78, 7).Value = "Poroto verde"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 60
$ws.Cells.Item(78, 11).Value = 25000
$ws.Cells.Item(78, 12).Value = 26000
$ws.Cells.Item(78, 13).Value = 25500
$ws.Cells.Item(78, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(78, 15).Value = "Región del Maule"
$ws.Cells.Item(78, 16).Value = 1020
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = "Hortaliza"
